# Update crypto price/volume figures per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.810.15"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "'2.094.82"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'228.84"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "'61.43"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").Value = "'0.0844"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "'15.35"
$ws.Range("E12").Value = "  +4.80%  "
$ws.Range("D13").Value = "'2.405.30"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'22.11"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "'0.810"
$ws.Range("E15").Value = "  +4.62%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "'2.092.63"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'38.758.67"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "'71.92"
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "'228.13"
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "'171.37"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("E29").Value = "  +3.70%  "
$ws.Range("D30").Value = "'19.31"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").Value = "'6.52"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'3.59"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'18.13"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("D42").Value = "'101.30"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'1.535.78"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("E47").Value = "  +6.00%  "
$ws.Range("D48").Value = "'4.11"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").Value = "'2.293.16"
